# WISCONSIN_2015.xlsx cleanup
#  1. Rename the header row to short machine-friendly column names.
#  2. Title-case the lowercase Spanish connector words ("de", "del", "el",
#     "la", "las", "los", "y") inside the state (A) and municipality (B)
#     text columns, e.g. "Pabellón de Arteaga" -> "Pabellón De Arteaga".
#  3. Clean up the stray carriage-return artifact in the "Estado de México"
#     state-name cell (A181) and title-case it the same way.
#  4. Drop the trailing footnote/source rows (1390-1395) that sit below the
#     data table, and shrink the sheet dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row -> machine-friendly names -------------------------------
$ws.Range("A1").Value2 = "mx_state"
$ws.Range("B1").Value2 = "mx_municipality"
$ws.Range("C1").Value2 = "n_matriculas"
$ws.Range("D1").Value2 = "pct_matriculas"

# --- 2. Title-case connector words in the state/municipality columns ------
$dataColumns = @($ws.Range("A2:A1389"), $ws.Range("B2:B1389"))
foreach ($col in $dataColumns) {
    $col.Replace(" de ", " De ", -4160)  | Out-Null
    $col.Replace(" del ", " Del ", -4160) | Out-Null
    $col.Replace(" el ", " El ", -4160)  | Out-Null
    $col.Replace(" la ", " La ", -4160)  | Out-Null
    $col.Replace(" las ", " Las ", -4160) | Out-Null
    $col.Replace(" los ", " Los ", -4160) | Out-Null
    $col.Replace(" y ", " Y ", -4160)   | Out-Null
}

# --- 3. Fix the "Estado de México" cell (had a trailing CRLF artifact) ----
$ws.Range("A181").Value2 = "Estado De México"

# --- 4. Remove the trailing footnote rows and shrink the used range -------
$ws.Range("A1390:A1395").EntireRow.Delete() | Out-Null
